$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.285.60'
$ws.Range('E2').Value = '  -0.75%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.575.46'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.03'
$ws.Range('E5').Value = '  -0.10%  '
$ws.Range('E6').Value = '  -0.36%  '
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.48'
$ws.Range('E8').Value = '  -3.49%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '23.79'
$ws.Range('E9').Value = '  -0.46%  '
$ws.Range('E10').Value = '  -0.69%  '
$ws.Range('E11').Value = '  -0.56%  '
$ws.Range('E12').Value = '  +1.74%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.800.99'
$ws.Range('E13').Value = '  +0.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.571.47'
$ws.Range('E14').Value = '  -0.37%  '
$ws.Range('E15').Value = '  -0.67%  '
$ws.Range('E16').Value = '  -0.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '28.319.58'
$ws.Range('E17').Value = '  -0.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.69'
$ws.Range('E18').Value = '  -1.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '230.53'
$ws.Range('E19').Value = '  +0.69%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.43'
$ws.Range('E20').Value = '  +0.81%  '
$ws.Range('E21').Value = '  -1.00%  '
$ws.Range('E22').Value = '  +0.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.93'
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.03'
$ws.Range('E24').Value = '  -1.40%  '
$ws.Range('E25').Value = '  +3.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.56'
$ws.Range('E26').Value = '  +0.38%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.98'
$ws.Range('E27').Value = '  -0.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.37'
$ws.Range('E28').Value = '  -1.36%  '
$ws.Range('E29').Value = '  -0.88%  '
$ws.Range('E30').Value = '  +0.32%  '
$ws.Range('E31').Value = '  +3.36%  '
$ws.Range('E32').Value = '  -3.61%  '
$ws.Range('E33').Value = '  -0.29%  '
$ws.Range('E34').Value = '  -1.40%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.387.55'
$ws.Range('E35').Value = '  -0.27%  '
$ws.Range('E36').Value = '  +6.54%  '
$ws.Range('E37').Value = '  -2.74%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.36'
$ws.Range('E38').Value = '  +0.05%  '
$ws.Range('E39').Value = '  +2.97%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0163'
$ws.Range('E40').Value = '  -1.36%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.519'
$ws.Range('E41').Value = '  -2.17%  '
$ws.Range('E42').Value = '  +0.23%  '
$ws.Range('E43').Value = '  +2.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.785'
$ws.Range('E44').Value = '  -0.98%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0458'
$ws.Range('E45').Value = '  -1.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.39'
$ws.Range('E46').Value = '  -4.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.925'
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.712.95'
$ws.Range('E49').Value = '  +0.14%  '
$ws.Range('B50').Value = 'mCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.14'
$ws.Range('E50').Value = '  +0.49%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '85.26'
$ws.Range('E51').Value = '  -0.30%  '
